$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Update the version number paragraph from "v3.0" to "v2.5".
#    The text "v3.0" lives at the very start of paragraph 4
#    ("v" "3" "." "0"), so edit the two digit characters in place.
# ---------------------------------------------------------------------
$versionPara = $d.Paragraphs.Item(4)
$versionStart = $versionPara.Range.Start
$d.Range($versionStart + 1, $versionStart + 2).Text = "2"
$d.Range($versionStart + 3, $versionStart + 4).Text = "5"

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the end of that paragraph to the
#    very beginning of the document (right before the first run of the
#    title paragraph "Tutorial for Basic Jenkins Integration").
#
#    A direct Bookmarks.Add() call using a collapsed Range(0,0) snaps to
#    cover the whole first paragraph instead of staying collapsed, so we
#    work around it: temporarily split off an empty leading paragraph,
#    add the bookmark collapsed at its end (position 1, not position 0),
#    then remove the temporary paragraph break again. Removing that
#    break slides the bookmark back down to position 0 without ever
#    calling Add() directly on position 0.
# ---------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

$startRange = $d.Range(0, 0)
$startRange.InsertParagraphBefore()

$newBookmarkRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)

$d.Range(0, 1).Delete()
